$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 9906.81
$ws.Range("B7").Value = 9844.7900000000009
$ws.Range("C7").Value = 307.87
$ws.Range("D7").Value = 305.93
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = -0.63
$ws.Range("G7").Value = 42612.675358796296
$ws.Range("H7").Value = $true

# Match the date/time number format used by the rest of column G (style index 1)
# by copying formats from the cell above, rather than re-creating a numFmt.
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
